$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: set all target cell values (A1:C25) ---
$ws.Range("B1").Value = 'Ementa atual:'
$ws.Range("C1").Value = 'Ementa modificada (dados modificados em vermelho):'

$ws.Range("B2").Value = 'LOT2049'
$ws.Range("C2").Value = 'LOT2049'

$ws.Range("A3").Value = 'Nome:'
$ws.Range("B3").Value = ' Genética e Biotecnologia Vegetal'
$ws.Range("C3").Value = ' Genética e Biotecnologia Vegetal'

$ws.Range("A4").Value = 'Name:'
$ws.Range("B4").Value = 'Genetics and Plant Biotechnology'
$ws.Range("C4").Value = 'Genetics and Plant Biotechnology'

$ws.Range("A5").Value = 'Créditos-aula:'
$ws.Range("B5").Value = '3'
$ws.Range("C5").Value = '3'

$ws.Range("A6").Value = 'Créditos-trabalho'
$ws.Range("B6").Value = '0'
$ws.Range("C6").Value = '0'

$ws.Range("A7").Value = 'Carga horária:'
$ws.Range("B7").Value = '45 h'
$ws.Range("C7").Value = '45 h'

$ws.Range("A8").Value = 'Ativação:'
$ws.Range("B8").Value = '01/01/2019'
$ws.Range("C8").Value = '01/01/2019'

$ws.Range("A9").Value = 'Semestre ideal:'
$ws.Range("B9").Value = 'EB-7'
$ws.Range("C9").Value = 'EB-7'

$ws.Range("A10").Value = 'Objetivos:'
$ws.Range("B10").Value = '8711290 - Elisson Antônio da Costa Romanel'
$ws.Range("C10").Value = '8711290 - Elisson Antônio da Costa Romanel'

$ws.Range("A11").Value = 'Objectives:'

$ws.Range("A12").Value = 'Docentes responsáveis:'

$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'

$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("B14").Value = 'Plant Development; Plant Tissue Culture; Plant genome and gene regulation; Techniques for plant transformation; Biotechnology for improvement of yield and quality traits; Science, society and environmental impact of GM crops.'
$ws.Range("C14").Value = 'Plant Development; Plant Tissue Culture; Plant genome and gene regulation; Techniques for plant transformation; Biotechnology for improvement of yield and quality traits; Science, society and environmental impact of GM crops.'

$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '01/01/2019'
$ws.Range("C15").Value = '01/01/2019'

$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("B16").Value = '1. Introduction in plant biotechnology and agriculture2. Plant physiology and development3. Plant tissue culture4. Eukaryotic gene regulation5. Agrobacterium-mediated plant genetic transformation6. Biolistic and other non-Agrobacterium technologies of plant transformation7. Vectors for plant transformation8. Genes, traits of interest and Biotechnological strategies for engineering plants9. Risk and benefits associated with genetically modified (GM) plants10. Synthetic biology in plants'
$ws.Range("C16").Value = '1. Introduction in plant biotechnology and agriculture2. Plant physiology and development3. Plant tissue culture4. Eukaryotic gene regulation5. Agrobacterium-mediated plant genetic transformation6. Biolistic and other non-Agrobacterium technologies of plant transformation7. Vectors for plant transformation8. Genes, traits of interest and Biotechnological strategies for engineering plants9. Risk and benefits associated with genetically modified (GM) plants10. Synthetic biology in plants'

$ws.Range("A17").Value = 'Avaliação:'

$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '8711290 - Elisson Antônio da Costa Romanel'
$ws.Range("C18").Value = '8711290 - Elisson Antônio da Costa Romanel'

$ws.Range("A19").Value = 'Critério:'
$ws.Range("B19").Value = 'Notas - N distribuído no semestre. A composição das "N" fica critério do docente.'
$ws.Range("C19").Value = 'Notas - N distribuído no semestre. A composição das "N" fica critério do docente.'

$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("B20").Value = 'MF = (somatório de N)/número de N (adequando o valor de N, quando houver peso distinto para as Ns).'
$ws.Range("C20").Value = 'MF = (somatório de N)/número de N (adequando o valor de N, quando houver peso distinto para as Ns).'

$ws.Range("A21").Value = 'Bibliografia:'
$ws.Range("B21").Value = 'NF = (MF + PR)/2, onde PR é uma prova de recuperação.Prova de Recuperação (PR) para alunos com Média Final (MF) maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final (NF) igual ou maior do que 5,0.'
$ws.Range("C21").Value = 'NF = (MF + PR)/2, onde PR é uma prova de recuperação.Prova de Recuperação (PR) para alunos com Média Final (MF) maior ou igual a 3,0 e menor do que 5,0. Será considerado aprovado o aluno que tenha obtido Nota Final (NF) igual ou maior do que 5,0.'

$ws.Range("A22").Value = 'Requisitos:'

$ws.Range("B23").Value = 'LOT2008 -  Bioquímica II  (Requisito fraco)
'
$ws.Range("C23").Value = 'LOT2008 -  Bioquímica II  (Requisito fraco)
'

$ws.Range("B24").Value = 'LOT2040 -  Engenharia Genética  (Requisito fraco)
'
$ws.Range("C24").Value = 'LOT2040 -  Engenharia Genética  (Requisito fraco)
'

$ws.Range("B25").Value = 'LOT2053 -  Microbiologia  (Requisito fraco)
'
$ws.Range("C25").Value = 'LOT2053 -  Microbiologia  (Requisito fraco)
'

# --- Step 2: remove the now-obsolete row 26 (shrinks dimension to A1:C25) ---
$ws.Rows.Item(26).Delete()

# --- Step 3: fix up row heights to match target layout ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(23).RowHeight = 30
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(22).AutoFit()
